# Versión 2.2 multithreading activado
# Update "Nuevo Precio" (column F) values for several SKUs as per the
# price list refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-breaking space (U+00A0) is used between "$" and the amount,
# matching the existing convention in this sheet.
$nbsp = [char]0x00A0

$ws.Range("F3").Value  = "$" + $nbsp + "2.776,77"
$ws.Range("F4").Value  = "$" + $nbsp + "1.295,77"
$ws.Range("F7").Value  = "$" + $nbsp + "1.573,45"
$ws.Range("F8").Value  = "$" + $nbsp + "1.573,45"
$ws.Range("F9").Value  = "$" + $nbsp + "1.573,45"
$ws.Range("F10").Value = "$" + $nbsp + "1.203,20"
$ws.Range("F12").Value = "$" + $nbsp + "629,33"
$ws.Range("F14").Value = "$" + $nbsp + "2.036,26"
